# cost_14.xlsx — "data : case 1"
# Update the three recorded cost values on Sheet1 and narrow column A by one
# character unit, matching the author's commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data values -----------------------------------------------------
$ws.Range("A1").Value = 148.80210896933113
$ws.Range("B1").Value = 4.7124185050074212
$ws.Range("C1").Value = 0.59911937377690805

# --- column A width: 11.7109375 -> 10.7109375 (one character narrower)
# Excel's ColumnWidth is expressed in "characters of the Normal font" and is
# itself rounded to the nearest screen pixel, so we feed it the value whose
# round-trip lands closest to the target stored width.
$ws.Columns("A:A").ColumnWidth = 9.85
